$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.653906
$ws.Range("H2").Value = 1.961718
$ws.Range("I2").Value = 0.00670030715761011
$ws.Range("J2").Value = 0.00670030715761011
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.778439
$ws.Range("N2").Value = 11.335317
$ws.Range("O2").Value = 0.4252971528324392
$ws.Range("P2").Value = 0.4252971528324392
$ws.Range("Q2").Value = 2.470743932734
$ws.Range("R2").Value = 22.236695394606
$ws.Range("S2").Value = 0.002849621557234393
$ws.Range("T2").Value = 0.002849621557234393
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.653906
$ws.Range("H3").Value = 1.961718
$ws.Range("I3").Value = 0.00670030715761011
$ws.Range("J3").Value = 0.00670030715761011
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("N3").Value = 13.00021
$ws.Range("O3").Value = 0.4877633593505858
$ws.Range("P3").Value = 0.4877633593505858
$ws.Range("Q3").Value = 2.833638440086666
$ws.Range("R3").Value = 25.50274596078
$ws.Range("S3").Value = 0.003268164327876682
$ws.Range("T3").Value = 0.003268164327876682
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.653906
$ws.Range("H4").Value = 1.961718
$ws.Range("I4").Value = 0.00670030715761011
$ws.Range("J4").Value = 0.00670030715761011
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2909853333333334
$ws.Range("N4").Value = 0.8729560000000001
$ws.Range("O4").Value = 0.03275300561492853
$ws.Range("P4").Value = 0.03275300561492853
$ws.Range("Q4").Value = 0.1902770553786667
$ws.Range("R4").Value = 1.712493498408
$ws.Range("S4").Value = 0.0002194551979549498
$ws.Range("T4").Value = 0.0002194551979549498
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.653906
$ws.Range("H5").Value = 1.961718
$ws.Range("I5").Value = 0.00670030715761011
$ws.Range("J5").Value = 0.00670030715761011
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4814053333333333
$ws.Range("N5").Value = 1.444216
$ws.Range("O5").Value = 0.0541864822020464
$ws.Range("P5").Value = 0.05418648220204641
$ws.Range("Q5").Value = 0.3147938358986667
$ws.Range("R5").Value = 2.833144523087999
$ws.Range("S5").Value = 0.0003630660745440843
$ws.Range("T5").Value = 0.0003630660745440844
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 96.11977900000001
$ws.Range("H6").Value = 288.359337
$ws.Range("I6").Value = 0.984900036429704
$ws.Range("J6").Value = 0.984900036429704
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.778439
$ws.Range("N6").Value = 11.335317
$ws.Range("O6").Value = 0.4252971528324392
$ws.Range("P6").Value = 0.4252971528324392
$ws.Range("Q6").Value = 363.182721644981
$ws.Range("R6").Value = 3268.644494804829
$ws.Range("S6").Value = 0.4188751813181188
$ws.Range("T6").Value = 0.4188751813181188
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 96.11977900000001
$ws.Range("H7").Value = 288.359337
$ws.Range("I7").Value = 0.984900036429704
$ws.Range("J7").Value = 0.984900036429704
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("N7").Value = 13.00021
$ws.Range("O7").Value = 0.4877633593505858
$ws.Range("P7").Value = 0.4877633593505858
$ws.Range("Q7").Value = 416.5257707178633
$ws.Range("R7").Value = 3748.73193646077
$ws.Range("S7").Value = 0.4803981503934668
$ws.Range("T7").Value = 0.4803981503934668
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 96.11977900000001
$ws.Range("H8").Value = 288.359337
$ws.Range("I8").Value = 0.984900036429704
$ws.Range("J8").Value = 0.984900036429704
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2909853333333334
$ws.Range("N8").Value = 0.8729560000000001
$ws.Range("O8").Value = 0.03275300561492853
$ws.Range("P8").Value = 0.03275300561492853
$ws.Range("Q8").Value = 27.96944593224134
$ws.Range("R8").Value = 251.7250133901721
$ws.Range("S8").Value = 0.03225843642332541
$ws.Range("T8").Value = 0.03225843642332541
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 96.11977900000001
$ws.Range("H9").Value = 288.359337
$ws.Range("I9").Value = 0.984900036429704
$ws.Range("J9").Value = 0.984900036429704
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.4814053333333333
$ws.Range("N9").Value = 1.444216
$ws.Range("O9").Value = 0.0541864822020464
$ws.Range("P9").Value = 0.05418648220204641
$ws.Range("Q9").Value = 46.27257424942133
$ws.Range("R9").Value = 416.453168244792
$ws.Range("S9").Value = 0.05336826829479301
$ws.Range("T9").Value = 0.05336826829479302
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2341223333333333
$ws.Range("H10").Value = 0.702367
$ws.Range("I10").Value = 0.002398955730318598
$ws.Range("J10").Value = 0.002398955730318598
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.778439
$ws.Range("N10").Value = 11.335317
$ws.Range("O10").Value = 0.4252971528324392
$ws.Range("P10").Value = 0.4252971528324392
$ws.Range("Q10").Value = 0.8846169550376667
$ws.Range("R10").Value = 7.961552595339
$ws.Range("S10").Value = 0.001020269041875565
$ws.Range("T10").Value = 0.001020269041875565
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2341223333333333
$ws.Range("H11").Value = 0.702367
$ws.Range("I11").Value = 0.002398955730318598
$ws.Range("J11").Value = 0.002398955730318598
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("N11").Value = 13.00021
$ws.Range("O11").Value = 0.4877633593505858
$ws.Range("P11").Value = 0.4877633593505858
$ws.Range("Q11").Value = 1.014546499674444
$ws.Range("R11").Value = 9.130918497069999
$ws.Range("S11").Value = 0.001170122705953538
$ws.Range("T11").Value = 0.001170122705953537
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.2341223333333333
$ws.Range("H12").Value = 0.702367
$ws.Range("I12").Value = 0.002398955730318598
$ws.Range("J12").Value = 0.002398955730318598
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.2909853333333334
$ws.Range("N12").Value = 0.8729560000000001
$ws.Range("O12").Value = 0.03275300561492853
$ws.Range("P12").Value = 0.03275300561492853
$ws.Range("Q12").Value = 0.06812616520577779
$ws.Range("R12").Value = 0.613135486852
$ws.Range("S12").Value = [double]"7.857301050509004E-05"
$ws.Range("T12").Value = [double]"7.857301050509002E-05"
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.2341223333333333
$ws.Range("H13").Value = 0.702367
$ws.Range("I13").Value = 0.002398955730318598
$ws.Range("J13").Value = 0.002398955730318598
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.4814053333333333
$ws.Range("N13").Value = 1.444216
$ws.Range("O13").Value = 0.0541864822020464
$ws.Range("P13").Value = 0.05418648220204641
$ws.Range("Q13").Value = 0.1127077399191111
$ws.Range("R13").Value = 1.014369659272
$ws.Range("S13").Value = 0.000129990971984406
$ws.Range("T13").Value = 0.000129990971984406
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5856290000000001
$ws.Range("H14").Value = 1.756887
$ws.Range("I14").Value = 0.00600070068236727
$ws.Range("J14").Value = 0.006000700682367269
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.778439
$ws.Range("N14").Value = 11.335317
$ws.Range("O14").Value = 0.4252971528324392
$ws.Range("P14").Value = 0.4252971528324392
$ws.Range("Q14").Value = 2.212763453131
$ws.Range("R14").Value = 19.914871078179
$ws.Range("S14").Value = 0.002552080915210475
$ws.Range("T14").Value = 0.002552080915210475
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5856290000000001
$ws.Range("H15").Value = 1.756887
$ws.Range("I15").Value = 0.00600070068236727
$ws.Range("J15").Value = 0.006000700682367269
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("N15").Value = 13.00021
$ws.Range("O15").Value = 0.4877633593505858
$ws.Range("P15").Value = 0.4877633593505858
$ws.Range("Q15").Value = 2.537766660696667
$ws.Range("R15").Value = 22.83989994627
$ws.Range("S15").Value = 0.002926921923288812
$ws.Range("T15").Value = 0.002926921923288812
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5856290000000001
$ws.Range("H16").Value = 1.756887
$ws.Range("I16").Value = 0.00600070068236727
$ws.Range("J16").Value = 0.006000700682367269
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.2909853333333334
$ws.Range("N16").Value = 0.8729560000000001
$ws.Range("O16").Value = 0.03275300561492853
$ws.Range("P16").Value = 0.03275300561492853
$ws.Range("Q16").Value = 0.1704094497746667
$ws.Range("R16").Value = 1.533685047972
$ws.Range("S16").Value = 0.0001965409831430807
$ws.Range("T16").Value = 0.0001965409831430806
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5856290000000001
$ws.Range("H17").Value = 1.756887
$ws.Range("I17").Value = 0.00600070068236727
$ws.Range("J17").Value = 0.006000700682367269
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.4814053333333333
$ws.Range("N17").Value = 1.444216
$ws.Range("O17").Value = 0.0541864822020464
$ws.Range("P17").Value = 0.05418648220204641
$ws.Range("Q17").Value = 0.2819249239546667
$ws.Range("R17").Value = 2.537324315592
$ws.Range("S17").Value = 0.0003251568607249018
$ws.Range("T17").Value = 0.0003251568607249018

Write-Host "Applied Dr Hou advice updates"